# "Rendu final du preTPI" - final corrections before submission:
#  - fix two typos in the Journal entries ("Rédacton" -> "Rédaction",
#    "Entretier" -> "Entretien")
#  - correct a logged duration (30 min -> 1h) on the last Journal entry
#  - make the Journal sheet the active sheet/selection again (the workbook
#    had been left with Totaux active / scrolled into the Journal data)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Journal")
$ws2 = $wb.Worksheets.Item("Totaux")

# Fix typos in the "Description" column of the Journal sheet.
$ws1.Range("E64").Value = "Rédaction du rapport de projet"
$ws1.Range("E72").Value = "Entretien avec la responsable du service de Qualité & Durabilité "

# Correct the logged time for the last entry (row 71): 30 min -> 1 h.
$ws1.Range("C71").Value = 1/24

# Restore the view: Journal active with E6 selected, Totaux no longer
# the active tab, with C26 selected there.
$ws2.Activate() | Out-Null
$ws2.Range("C26").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("E6").Select() | Out-Null
